# Applies the "borrar datos de cache al cerrar sesion" update to the
# resumen_movimiento_cuenta workbook:
#   - a handful of rows get newly-recorded "total_cobro" (J) / "total_fcp2" (G)
#     / "total_fcp2" (H) movements, which change the computed "saldo" (M)
#   - the "total_ne" (E) column is reformatted with a thousands-separator
#     ("Comma"/"Millares") number format
#   - a new totals row (104) is appended, summing E and M and showing the
#     variance between them in N
#   - the view is scrolled back to the top and a couple of columns are
#     resized to fit their (now wider) content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row-level value corrections (new collection/adjustment amounts posted
# against a handful of customer accounts, which shift their "saldo").
# ---------------------------------------------------------------------

# Row 3
$ws.Range("J3").Value = 98.9
$ws.Range("M3").Value = 16765.38

# Row 14
$ws.Range("E14").Value = 1105.77
$ws.Range("J14").Value = 106.02
$ws.Range("M14").Value = 999.75

# Row 39
$ws.Range("J39").Value = 1049.51
$ws.Range("M39").Value = -2135.35

# Row 54
$ws.Range("J54").Value = 306
$ws.Range("M54").Value = -998.72

# Row 66
$ws.Range("J66").Value = 2197.8
$ws.Range("M66").Value = 5187.95

# Row 68
$ws.Range("H68").Value = 99.61
$ws.Range("J68").Value = 416.77
$ws.Range("M68").Value = -8035.96

# Row 88
$ws.Range("E88").Value = 1091.79
$ws.Range("M88").Value = 1091.79

# Row 101
$ws.Range("G101").Value = 315
$ws.Range("M101").Value = 2808.28

# ---------------------------------------------------------------------
# Apply the thousands-separator ("Comma"/"Millares") number format to the
# whole total_ne (E) column of data.
# ---------------------------------------------------------------------
$ws.Range("E2:E103").NumberFormat = '_(* #,##0.00_);_(* \(#,##0.00\);_(* "-"??_);_(@_)'

# ---------------------------------------------------------------------
# New totals row: sums of E and M, plus the variance between them.
# ---------------------------------------------------------------------
$ws.Range("E104").Formula = "=SUM(E2:E103)"
$ws.Range("E104").NumberFormat = '_(* #,##0.00_);_(* \(#,##0.00\);_(* "-"??_);_(@_)'

$ws.Range("M104").Formula = "=SUM(M2:M103)"
$ws.Range("M104").NumberFormat = "[Blue]#,##0.00;[Red]\(#,##0.00\);[Black]#,##0.00"

$ws.Range("N104").Formula = "=M104-E104"
$ws.Range("N104").NumberFormat = "#,##0.00;[Red]#,##0.00"

$wb.Application.Calculate()

# ---------------------------------------------------------------------
# Column widths: widen E (now showing comma-separated totals) and bump M
# to fit its (slightly) wider formatted numbers.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 11
$ws.Columns.Item(13).ColumnWidth = 10.140625

# ---------------------------------------------------------------------
# Reset the view: scroll back to the top and select M1 (previously the
# sheet was left scrolled down to row 71 with M2:M103 selected).
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("M1").Select()
